$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2882.4
$ws.Range("I113").Value = 2347.1428
$ws.Range("J113").Value = 3350.75
$ws.Range("K113").Value = 2347.1428
$ws.Range("L113").Value = 3350.75
$ws.Range("M113").Value = 906.8571999999999
$ws.Range("N113").Value = -9858.75

$ws.Range("H129").Value = 1026.2264
$ws.Range("J129").Value = 1099.9783
$ws.Range("L129").Value = 3299.9349
$ws.Range("N129").Value = -13299.9349

$ws.Range("H137").Value = 1126.4375
$ws.Range("J137").Value = 2358.125
$ws.Range("L137").Value = 7074.375
$ws.Range("N137").Value = -12174.375

$ws.Range("H138").Value = 2344.2
$ws.Range("I138").Value = 1367.34
$ws.Range("J138").Value = 3972.3
$ws.Range("K138").Value = 4102.02
$ws.Range("L138").Value = 11916.9
$ws.Range("M138").Value = 1037.98
$ws.Range("N138").Value = -22196.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 201697.8
$ws.Range("I2").Value = 1996.3334
$ws.Range("J2").Value = 501250
$ws.Range("K2").Value = 1996.3334
$ws.Range("L2").Value = 501250
$ws.Range("M2").Value = -1883.3334
$ws.Range("N2").Value = -501476

$ws.Range("H74").Value = 866
$ws.Range("I74").Value = 1055.7407
$ws.Range("J74").Value = 622.0476
$ws.Range("K74").Value = 1055.7407
$ws.Range("L74").Value = 622.0476
$ws.Range("M74").Value = -181.7407000000001
$ws.Range("N74").Value = -2370.0476

$ws.Range("H77").Value = 866
$ws.Range("I77").Value = 1055.7407
$ws.Range("J77").Value = 622.0476
$ws.Range("K77").Value = 5278.703500000001
$ws.Range("L77").Value = 3110.238
$ws.Range("M77").Value = -910.7035000000005
$ws.Range("N77").Value = -11846.238

$ws.Range("H110").Value = 1513.0769
$ws.Range("I110").Value = 1464.7
$ws.Range("J110").Value = 1674.3334
$ws.Range("K110").Value = 1464.7
$ws.Range("L110").Value = 1674.3334
$ws.Range("M110").Value = 580.3
$ws.Range("N110").Value = -5764.3334

$ws.Range("H116").Value = 201697.8
$ws.Range("I116").Value = 1996.3334
$ws.Range("J116").Value = 501250
$ws.Range("K116").Value = 1996.3334
$ws.Range("L116").Value = 501250
$ws.Range("M116").Value = 297.6666
$ws.Range("N116").Value = -505838

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 201697.8
$ws.Range("I3").Value = 1996.3334
$ws.Range("J3").Value = 501250
$ws.Range("K3").Value = 1996.3334
$ws.Range("L3").Value = 501250
$ws.Range("M3").Value = -1882.3334
$ws.Range("N3").Value = -501478

$ws.Range("H105").Value = 2559.4783
$ws.Range("I105").Value = 2426.611
$ws.Range("J105").Value = 3037.8
$ws.Range("K105").Value = 2426.611
$ws.Range("L105").Value = 3037.8
$ws.Range("M105").Value = -679.6109999999999
$ws.Range("N105").Value = -6531.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 4613
$ws.Range("I22").Value = 5761.5557
$ws.Range("J22").Value = 478.2
$ws.Range("K22").Value = 5761.5557
$ws.Range("L22").Value = 478.2
$ws.Range("M22").Value = -5411.5557
$ws.Range("N22").Value = -1178.2

$ws.Range("H31").Value = 2099.5967
$ws.Range("I31").Value = 1609.738
$ws.Range("J31").Value = 3128.3
$ws.Range("K31").Value = 1609.738
$ws.Range("L31").Value = 3128.3
$ws.Range("M31").Value = -1314.738
$ws.Range("N31").Value = -3718.3

$ws.Range("H34").Value = 2099.5967
$ws.Range("I34").Value = 1609.738
$ws.Range("J34").Value = 3128.3
$ws.Range("K34").Value = 1609.738
$ws.Range("L34").Value = 3128.3
$ws.Range("M34").Value = -1407.738
$ws.Range("N34").Value = -3532.3

$ws.Range("H99").Value = 2814.4736
$ws.Range("I99").Value = 2658.3333
$ws.Range("J99").Value = 3400
$ws.Range("K99").Value = 2658.3333
$ws.Range("L99").Value = 3400
$ws.Range("M99").Value = -1160.3333
$ws.Range("N99").Value = -6396

$ws.Range("H126").Value = 2814.4736
$ws.Range("I126").Value = 2658.3333
$ws.Range("J126").Value = 3400
$ws.Range("K126").Value = 7974.999899999999
$ws.Range("L126").Value = 10200
$ws.Range("M126").Value = -5504.999899999999
$ws.Range("N126").Value = -15140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 602.125
$ws.Range("I113").Value = 500
$ws.Range("J113").Value = 733.4286
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 2200.2858
$ws.Range("M113").Value = 670
$ws.Range("N113").Value = -6540.2858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6811
$ws.Range("I70").Value = 6080
$ws.Range("J70").Value = 7623.222
$ws.Range("K70").Value = 6080
$ws.Range("L70").Value = 7623.222
$ws.Range("M70").Value = -5810
$ws.Range("N70").Value = -8163.222

$ws.Range("H73").Value = 6811
$ws.Range("I73").Value = 6080
$ws.Range("J73").Value = 7623.222
$ws.Range("K73").Value = 6080
$ws.Range("L73").Value = 7623.222
$ws.Range("M73").Value = -5144
$ws.Range("N73").Value = -9495.222

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 857.1429000000001
$ws.Range("I22").Value = 825
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 825
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -530
$ws.Range("N22").Value = -1490

$ws.Range("H27").Value = 857.1429000000001
$ws.Range("I27").Value = 825
$ws.Range("J27").Value = 900
$ws.Range("K27").Value = 825
$ws.Range("L27").Value = 900
$ws.Range("M27").Value = -718
$ws.Range("N27").Value = -1114

$ws.Range("H61").Value = 38167.332
$ws.Range("I61").Value = 45000.8
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 45000.8
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -44798.8
$ws.Range("N61").Value = -4404

$ws.Range("H113").Value = 38167.332
$ws.Range("I113").Value = 45000.8
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 45000.8
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -42830.8
$ws.Range("N113").Value = -8340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 41096.062
$ws.Range("I81").Value = 38839.63
$ws.Range("J81").Value = 51250
$ws.Range("K81").Value = 77679.25999999999
$ws.Range("L81").Value = 102500
$ws.Range("M81").Value = -76618.25999999999
$ws.Range("N81").Value = -104622

$ws.Range("H84").Value = 41096.062
$ws.Range("I84").Value = 38839.63
$ws.Range("J84").Value = 51250
$ws.Range("K84").Value = 388396.3
$ws.Range("L84").Value = 512500
$ws.Range("M84").Value = -383092.3
$ws.Range("N84").Value = -523108

$ws.Range("H107").Value = 506.4737
$ws.Range("I107").Value = 543.1667
$ws.Range("J107").Value = 443.57144
$ws.Range("K107").Value = 1629.5001
$ws.Range("L107").Value = 1330.71432
$ws.Range("M107").Value = 290.4999
$ws.Range("N107").Value = -5170.71432

$ws.Range("H122").Value = 11576568
$ws.Range("I122").Value = 17859180
$ws.Range("K122").Value = 53577540
$ws.Range("M122").Value = -53575090

$ws.Range("H132").Value = 1044.4073
$ws.Range("I132").Value = 673.13336
$ws.Range("J132").Value = 2900.7778
$ws.Range("K132").Value = 2019.40008
$ws.Range("L132").Value = 8702.3334
$ws.Range("M132").Value = 510.5999199999999
$ws.Range("N132").Value = -13762.3334

$ws.Range("H136").Value = 1460.0555
$ws.Range("I136").Value = 1555.2667
$ws.Range("J136").Value = 984
$ws.Range("K136").Value = 4665.800099999999
$ws.Range("L136").Value = 2952
$ws.Range("M136").Value = -2115.800099999999
$ws.Range("N136").Value = -8052
